$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[KNIME] Loop 실습(1)"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/KNIME-Loop-%EC%8B%A4%EC%8A%B5"

$ws.Range("D32").Value = "추천 메트릭 :: Precision@k, Recall@k, MAP, MRR, NDCG , AP, F1-Score, Coverage, Diversity, Novelty"
$ws.Range("E32").Value = "https://dodonam.tistory.com/493"

$ws.Range("D51").Value = "PDF에 적용된 폰트 확인하는 방법 (Adobe Acrobat Reader)"
$ws.Range("E51").Value = "https://bskyvision.com/entry/PDF%EC%97%90-%EC%A0%81%EC%9A%A9%EB%90%9C-%ED%8F%B0%ED%8A%B8-%ED%99%95%EC%9D%B8%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95-Adobe-Acrobat-Reader"
